# Apply the changes described by the commit:
#   "completing the square added / new study guide, questions, answers :)"
#
# 1) Insert a new "Summary" paragraph (style AbstractTitle) right after the
#    Author paragraph ("Ifan Howells-Baines, Mark Toner") and before the
#    existing Abstract paragraph ("tbc").
# 2) Character-style tweaks in styles.xml (syntax-highlighting tokens):
#      - KeywordTok / ControlFlowTok gain bold.
#      - DocumentationTok / CommentVarTok / WarningTok keep their italic but
#        it is re-asserted so it sits first in the run properties.

$d = $word.ActiveDocument

# --- 1. Insert the "Summary" / AbstractTitle paragraph -------------------

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*Toner*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $authorPara = $d.Paragraphs.Item($targetIndex)
    $authorPara.Range.InsertParagraphAfter()

    $summaryPara = $d.Paragraphs.Item($targetIndex + 1)
    $summaryPara.Range.Text = "Summary"
    $summaryPara.Style = "AbstractTitle"
}

# --- 2. Re-order / add run-property flags on the Pandoc Tok styles -------

$d.Styles.Item("KeywordTok").Font.Bold = $true
$d.Styles.Item("ControlFlowTok").Font.Bold = $true

$d.Styles.Item("DocumentationTok").Font.Italic = $true
$d.Styles.Item("CommentVarTok").Font.Italic = $true
$d.Styles.Item("WarningTok").Font.Italic = $true
